# The "H 72" data row (row 2) was removed from the sheet; all subsequent
# rows shift up by one and the used range shrinks from A1:F63 to A1:F62.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Delete()
